$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price-report row is inserted at row 190 (the data is sorted
# newest-first by date), pushing the former rows 190-292 down to 191-293.
$ws.Rows.Item(190).Insert("xlShiftDown")

# Populate the newly inserted row with this week's Cilantro price record.
$ws.Range("A190").Value = 10
$ws.Range("B190").Value = "Vega Modelo de Temuco"
$ws.Range("C190").Value = "La Araucanía"
$ws.Range("D190").Value = 44572
$ws.Range("E190").Value = 9
$ws.Range("F190").Value = 100112040
$ws.Range("G190").Value = "Cilantro"
$ws.Range("H190").Value = "Sin especificar"
$ws.Range("I190").Value = "Primera"
$ws.Range("J190").Value = 75
$ws.Range("K190").Value = 8000
$ws.Range("L190").Value = 8000
$ws.Range("M190").Value = 8000
$ws.Range("N190").Value = "`$/docena de atados (2 kilos)"
$ws.Range("O190").Value = "Provincia de Cautín"
$ws.Range("P190").Value = 4000
$ws.Range("Q190").Value = 2
$ws.Range("R190").Value = "Hortaliza"
